{"js": "// Update the date line and the 25 multiplication problems in the table.\nconst body = context.document.body;\n\n// --- 1. Date paragraph (first paragraph in the body) ---\nconst firstPara = body.paragraphs.getFirst();\nfirstPara.load(\"text\");\nawait context.sync();\n\nif (firstPara.text.indexOf(\"2024-06-13 Thursday\") !== -1) {\n  const dateRange = firstPara.getRange();\n  dateRange.insertText(\"2024-06-14 Friday\", \"Replace\");\n}\n\n// --- 2. Table of multiplication problems ---\nconst table = body.tables.getFirst();\n\n// Each inner array is one non-blank table row (by row index within the\n// table), holding the old->new text for each of the 5 columns.\nconst rowEdits = [\n  { row: 0, cells: [\"510\u00d77=\", \"644\u00d73=\", \"773\u00d73=\", \"524\u00d76=\", \"224\u00d78=\"] },\n  { row: 4, cells: [\"365\u00d76=\", \"527\u00d77=\", \"793\u00d77=\", \"661\u00d79=\", \"407\u00d76=\"] },\n  { row: 9, cells: [\"388\u00d73=\", \"268\u00d72=\", \"259\u00d76=\", \"453\u00d75=\", \"938\u00d78=\"] },\n  { row: 14, cells: [\"490\u00d77=\", \"350\u00d79=\", \"661\u00d79=\", \"907\u00d72=\", \"420\u00d74=\"] },\n  { row: 19, cells: [\"263\u00d76=\", \"949\u00d73=\", \"463\u00d78=\", \"649\u00d78=\", \"781\u00d78=\"] },\n];\n\nfor (const { row, cells } of rowEdits) {\n  for (let col = 0; col < cells.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = cells[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Date paragraph ---\n$found = $d.Content.Find.Execute(\"2024-06-13 Thursday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-06-14 Friday\", 2)\n\n# --- 2. Table of multiplication problems ---\n$t = $d.Tables.Item(1)\n\n$rowEdits = @{\n    1 = @(\"510\u00d77=\", \"644\u00d73=\", \"773\u00d73=\", \"524\u00d76=\", \"224\u00d78=\")\n    5 = @(\"365\u00d76=\", \"527\u00d77=\", \"793\u00d77=\", \"661\u00d79=\", \"407\u00d76=\")\n    10 = @(\"388\u00d73=\", \"268\u00d72=\", \"259\u00d76=\", \"453\u00d75=\", \"938\u00d78=\")\n    15 = @(\"490\u00d77=\", \"350\u00d79=\", \"661\u00d79=\", \"907\u00d72=\", \"420\u00d74=\")\n    20 = @(\"263\u00d76=\", \"949\u00d73=\", \"463\u00d78=\", \"649\u00d78=\", \"781\u00d78=\")\n}\n\nforeach ($row in $rowEdits.Keys) {\n    $vals = $rowEdits[$row]\n    for ($col = 1; $col -le $vals.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $vals[$col - 1]\n    }\n}\n"}
